$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset has shifted forward by 13 days (new reporting period),
# and the "Notified Production (MW)" forecast values (column B) have been revised
# for the NRG, PCSun and Ulmeni models.
$n = 96
$data = New-Object 'object[,]' $n,2

$data[0,0] = 46069.01041666666
$data[0,1] = 0
$data[1,0] = 46069.02083333334
$data[1,1] = 0
$data[2,0] = 46069.03125
$data[2,1] = 0
$data[3,0] = 46069.04166666666
$data[3,1] = 0
$data[4,0] = 46069.05208333334
$data[4,1] = 0.39
$data[5,0] = 46069.0625
$data[5,1] = 0
$data[6,0] = 46069.07291666666
$data[6,1] = 0
$data[7,0] = 46069.08333333334
$data[7,1] = 0
$data[8,0] = 46069.09375
$data[8,1] = 0
$data[9,0] = 46069.10416666666
$data[9,1] = 0
$data[10,0] = 46069.11458333334
$data[10,1] = 0
$data[11,0] = 46069.125
$data[11,1] = 0
$data[12,0] = 46069.13541666666
$data[12,1] = 0.55
$data[13,0] = 46069.14583333334
$data[13,1] = 0
$data[14,0] = 46069.15625
$data[14,1] = 0
$data[15,0] = 46069.16666666666
$data[15,1] = 0
$data[16,0] = 46069.17708333334
$data[16,1] = 0
$data[17,0] = 46069.1875
$data[17,1] = 0
$data[18,0] = 46069.19791666666
$data[18,1] = 0
$data[19,0] = 46069.20833333334
$data[19,1] = 0
$data[20,0] = 46069.21875
$data[20,1] = 0.6860000000000001
$data[21,0] = 46069.22916666666
$data[21,1] = 0.6870000000000001
$data[22,0] = 46069.23958333334
$data[22,1] = 0.711
$data[23,0] = 46069.25
$data[23,1] = 0.897
$data[24,0] = 46069.26041666666
$data[24,1] = 7.216
$data[25,0] = 46069.27083333334
$data[25,1] = 14.159
$data[26,0] = 46069.28125
$data[26,1] = 27.052
$data[27,0] = 46069.29166666666
$data[27,1] = 57.064
$data[28,0] = 46069.30208333334
$data[28,1] = 157.776
$data[29,0] = 46069.3125
$data[29,1] = 227.167
$data[30,0] = 46069.32291666666
$data[30,1] = 309.307
$data[31,0] = 46069.33333333334
$data[31,1] = 436.261
$data[32,0] = 46069.34375
$data[32,1] = 612.042
$data[33,0] = 46069.35416666666
$data[33,1] = 719.227
$data[34,0] = 46069.36458333334
$data[34,1] = 821.134
$data[35,0] = 46069.375
$data[35,1] = 915.16
$data[36,0] = 46069.38541666666
$data[36,1] = 1094.758
$data[37,0] = 46069.39583333334
$data[37,1] = 1169.21
$data[38,0] = 46069.40625
$data[38,1] = 1259.706
$data[39,0] = 46069.41666666666
$data[39,1] = 1328.105
$data[40,0] = 46069.42708333334
$data[40,1] = 1394.833
$data[41,0] = 46069.4375
$data[41,1] = 1439.283
$data[42,0] = 46069.44791666666
$data[42,1] = 1481.957
$data[43,0] = 46069.45833333334
$data[43,1] = 1511.195
$data[44,0] = 46069.46875
$data[44,1] = 1525.301
$data[45,0] = 46069.47916666666
$data[45,1] = 1541.198
$data[46,0] = 46069.48958333334
$data[46,1] = 1542.049
$data[47,0] = 46069.5
$data[47,1] = 1534.072
$data[48,0] = 46069.51041666666
$data[48,1] = 1489.419
$data[49,0] = 46069.52083333334
$data[49,1] = 1464.154
$data[50,0] = 46069.53125
$data[50,1] = 1424.968
$data[51,0] = 46069.54166666666
$data[51,1] = 1382.483
$data[52,0] = 46069.55208333334
$data[52,1] = 1290.453
$data[53,0] = 46069.5625
$data[53,1] = 1231.114
$data[54,0] = 46069.57291666666
$data[54,1] = 1171.393
$data[55,0] = 46069.58333333334
$data[55,1] = 1093.078
$data[56,0] = 46069.59375
$data[56,1] = 957.067
$data[57,0] = 46069.60416666666
$data[57,1] = 877.706
$data[58,0] = 46069.61458333334
$data[58,1] = 786.806
$data[59,0] = 46069.625
$data[59,1] = 695.623
$data[60,0] = 46069.63541666666
$data[60,1] = 529.395
$data[61,0] = 46069.64583333334
$data[61,1] = 431.763
$data[62,0] = 46069.65625
$data[62,1] = 341.098
$data[63,0] = 46069.66666666666
$data[63,1] = 261.715
$data[64,0] = 46069.67708333334
$data[64,1] = 125.296
$data[65,0] = 46069.6875
$data[65,1] = 79.922
$data[66,0] = 46069.69791666666
$data[66,1] = 52.545
$data[67,0] = 46069.70833333334
$data[67,1] = 27.497
$data[68,0] = 46069.71875
$data[68,1] = 10.098
$data[69,0] = 46069.72916666666
$data[69,1] = 6.993
$data[70,0] = 46069.73958333334
$data[70,1] = 5.998
$data[71,0] = 46069.75
$data[71,1] = 5.605
$data[72,0] = 46069.76041666666
$data[72,1] = 10.55
$data[73,0] = 46069.77083333334
$data[73,1] = 10.65
$data[74,0] = 46069.78125
$data[74,1] = 0
$data[75,0] = 46069.79166666666
$data[75,1] = 0
$data[76,0] = 46069.80208333334
$data[76,1] = 10.49
$data[77,0] = 46069.8125
$data[77,1] = 8.49
$data[78,0] = 46069.82291666666
$data[78,1] = 12.49
$data[79,0] = 46069.83333333334
$data[79,1] = 7.29
$data[80,0] = 46069.84375
$data[80,1] = 6.55
$data[81,0] = 46069.85416666666
$data[81,1] = 4.55
$data[82,0] = 46069.86458333334
$data[82,1] = 0
$data[83,0] = 46069.875
$data[83,1] = 2.55
$data[84,0] = 46069.88541666666
$data[84,1] = 4.55
$data[85,0] = 46069.89583333334
$data[85,1] = 2.55
$data[86,0] = 46069.90625
$data[86,1] = 0
$data[87,0] = 46069.91666666666
$data[87,1] = 0
$data[88,0] = 46069.92708333334
$data[88,1] = 0
$data[89,0] = 46069.9375
$data[89,1] = 0.55
$data[90,0] = 46069.94791666666
$data[90,1] = 0
$data[91,0] = 46069.95833333334
$data[91,1] = 0
$data[92,0] = 46069.96875
$data[92,1] = 0
$data[93,0] = 46069.97916666666
$data[93,1] = 0
$data[94,0] = 46069.98958333334
$data[94,1] = 0
$data[95,0] = 46070
$data[95,1] = 0

$range = $ws.Range("A2:B97")
$range.Value2 = $data

Write-Host "Updated A2:B97 with new timestamps (+13 days) and revised production forecast values."
